$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 325 and 326 (columns B and D changed, C unchanged) ---
$ws.Cells.Item(325, 2).Value = 7710432000000
$ws.Cells.Item(325, 4).Value = 255905476269.4988

$ws.Cells.Item(326, 2).Value = 7822810000000
$ws.Cells.Item(326, 4).Value = 255814584695.8797

# --- Append new rows 327-329, matching column A's existing date style (s="2") ---
$ws.Cells.Item(326, 1).Copy()
$ws.Cells.Item(327, 1).PasteSpecial(-4122)
$ws.Cells.Item(328, 1).PasteSpecial(-4122)
$ws.Cells.Item(329, 1).PasteSpecial(-4122)

$ws.Cells.Item(327, 1).Value = 44986
$ws.Cells.Item(327, 2).Value = 7965088000000
$ws.Cells.Item(327, 3).Value = 0.03252032520325204
$ws.Cells.Item(327, 4).Value = 259027252032.5204

$ws.Cells.Item(328, 1).Value = 45017
$ws.Cells.Item(328, 2).Value = 8069151000000
$ws.Cells.Item(328, 3).Value = 0.03241491085899514
$ws.Cells.Item(328, 4).Value = 261560810372.7715

$ws.Cells.Item(329, 1).Value = 45047
$ws.Cells.Item(329, 2).Value = 8140535000000
$ws.Cells.Item(329, 3).Value = 0.03241491085899514
$ws.Cells.Item(329, 4).Value = 263874716369.53
